$d = $word.ActiveDocument

# The paragraph containing the sentence we need to edit/re-run-split.
$paraIndex = 25

# --- Step 1: fix the text content -------------------------------------
# "princip S, te se"   -> "princip S te se"   (drop the comma)
# "slažu s navedenim"  -> "slažu sa navedenim" (s -> sa)
$r1 = $d.Paragraphs($paraIndex).Range
$null = $r1.Find.Execute("princip S, te se", $true, $false, $false, $false, `
                          $false, $true, 1, $false, "princip S te se", 2)

$r2 = $d.Paragraphs($paraIndex).Range
$null = $r2.Find.Execute("slažu s navedenim", $true, $false, $false, $false, `
                          $false, $true, 1, $false, "slažu sa navedenim", 2)

# --- Step 2: split the single run into five runs -----------------------
# Piece boundaries, relative to the (now corrected) paragraph text:
#   [0,89)    Clan tima Belma ... Ehlimana izrazila
#   [89,120)   je zelju da obradjuje princip S
#   [120,295)  te se izaslo ... da se slazu s
#   [295,296)  a
#   [296,308)  navedenim).
$p = $d.Paragraphs($paraIndex)
$pStart = $p.Range.Start
$pEnd = $p.Range.End

$b1 = $pStart + 89
$b2 = $pStart + 120
$b3 = $pStart + 295
$b4 = $pStart + 296

$d.Range($pStart, $b1).Font.Name = "Book Antiqua"
$d.Range($b1, $b2).Font.Name = "Book Antiqua"
$d.Range($b2, $b3).Font.Name = "Book Antiqua"
$d.Range($b3, $b4).Font.Name = "Book Antiqua"
$d.Range($b4, $pEnd - 1).Font.Name = "Book Antiqua"

# Restore the complex-script font name (w:cs) that Font.Name alone cannot
# set; this call affects every run in the paragraph, which is fine here
# since they all need the same "Times New Roman" cs font.
$d.Paragraphs($paraIndex).Range.Font.NameBi = "Times New Roman"

# --- Step 3: move the _GoBack bookmark --------------------------------
# Word keeps a single "_GoBack" bookmark tracking the last edit location;
# adding it at the new position removes it from wherever it used to be.
$bmRange = $d.Range($b4, $b4)
$d.Bookmarks.Add("_GoBack", $bmRange)
